$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    3  = @(0.3464964993005633, 0.3375848360084654, 3.082599426703578, 0.4998867070740569, 1, 4.266567469086664)
    4  = @(1.505614041169197, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245)
    5  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 12.0302756157461)
    6  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    7  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    8  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 3.594575437922795)
    9  = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538)
    10 = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 0, 2.27892381503245)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i  # B=2 .. G=7
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
